$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before the existing "F" (description) column, shifting
# it to "G". Excel copies the formatting of the column to the left (style 1),
# matching the source workbook's uniform cell style.
$ws.Columns("F:F").Insert()

# Make sure the new column uses the same "text" number format as the rest of
# the table (style index 1 == numFmtId 49 "@") so every cell in the column
# round-trips with s="1" once populated, matching its siblings.
$ws.Range("F1:F6").NumberFormat = "@"

# New column width (13.525 in the saved file). Excel's ColumnWidth (character
# units) gets rounded to whole pixels (Maximum Digit Width = 7) before being
# stored, so 12.81 is the input that lands on the closest achievable pixel
# bucket to the target width.
$ws.Range("F1").ColumnWidth = 12.81

# Header rows for the new "TargetCount" column.
$ws.Range("F1").Value = "目标数量"
$ws.Range("F2").Value = "TargetCount"
$ws.Range("F3").Value = "int"

# Data rows: every record gets a target count of 1.
$ws.Range("F4").Value = "1"
$ws.Range("F5").Value = "1"
$ws.Range("F6").Value = "1"

# Restore the selection state recorded in the saved workbook.
$ws.Range("H23").Select()
